# Append 7 new "Bag" data rows (rows 20-26) to WorkSheet 1, matching the
# existing table layout (columns A:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, timestamp serial (col A), ElapsedMs (col C)
# Remaining columns D:M are constant across all rows in this table:
# D=17 E=2 F=0 G=1 H=0 I=100 J=2 K=0 L=100 M=0
$newRows = @(
    @{ Row = 20; A = 42602.00953703704;  C = 46 },
    @{ Row = 21; A = 42602.014432870368; C = 51 },
    @{ Row = 22; A = 42602.015277777777; C = 15 },
    @{ Row = 23; A = 42602.495196759257; C = 47 },
    @{ Row = 24; A = 42602.495486111111; C = 18 },
    @{ Row = 25; A = 42602.495729166665; C = 15 },
    @{ Row = 26; A = 42602.495798611111; C = 13 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = "Bag"
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = 17
    $ws.Range("E$r").Value = 2
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 1
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = 100
    $ws.Range("J$r").Value = 2
    $ws.Range("K$r").Value = 0
    $ws.Range("L$r").Value = 100
    $ws.Range("M$r").Value = 0
}
